$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INTERNAL ISSUE")
$ws.Range("C10:C3062").NumberFormat = "mm/dd/yy;@"
